# Fruta / hortaliza, semanal
# Insert two new weekly records at the top of the data table (rows 67-68),
# pushing the existing records down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 67 (shifts old rows 67.. down to 69..)
$ws.Rows.Item(67).Insert()
$ws.Rows.Item(67).Insert()

# New row 67
$ws.Cells.Item(67, 1).Value = 9
$ws.Cells.Item(67, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(67, 3).Value = "Metropolitana"
$ws.Cells.Item(67, 4).Value = 44588
$ws.Cells.Item(67, 5).Value = 13
$ws.Cells.Item(67, 6).Value = "Fruta"
$ws.Cells.Item(67, 7).Value = 100101
$ws.Cells.Item(67, 8).Value = "Berries"
$ws.Cells.Item(67, 9).Value = 100101004
$ws.Cells.Item(67, 10).Value = "Frambuesa"
$ws.Cells.Item(67, 11).Value = "Sin especificar"
$ws.Cells.Item(67, 12).Value = "Especial"
$ws.Cells.Item(67, 13).Value = 310
$ws.Cells.Item(67, 14).Value = 8000
$ws.Cells.Item(67, 15).Value = 8000
$ws.Cells.Item(67, 16).Value = 8000
$ws.Cells.Item(67, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(67, 18).Value = "Provincia de Linares"
$ws.Cells.Item(67, 19).Value = 4000
$ws.Cells.Item(67, 20).Value = 2

# New row 68
$ws.Cells.Item(68, 1).Value = 9
$ws.Cells.Item(68, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(68, 3).Value = "Metropolitana"
$ws.Cells.Item(68, 4).Value = 44588
$ws.Cells.Item(68, 5).Value = 13
$ws.Cells.Item(68, 6).Value = "Fruta"
$ws.Cells.Item(68, 7).Value = 100101
$ws.Cells.Item(68, 8).Value = "Berries"
$ws.Cells.Item(68, 9).Value = 100101004
$ws.Cells.Item(68, 10).Value = "Frambuesa"
$ws.Cells.Item(68, 11).Value = "Sin especificar"
$ws.Cells.Item(68, 12).Value = "Primera"
$ws.Cells.Item(68, 13).Value = 350
$ws.Cells.Item(68, 14).Value = 7000
$ws.Cells.Item(68, 15).Value = 7000
$ws.Cells.Item(68, 16).Value = 7000
$ws.Cells.Item(68, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(68, 18).Value = "Provincia de Linares"
$ws.Cells.Item(68, 19).Value = 3500
$ws.Cells.Item(68, 20).Value = 2
